# Change the "unit" value for all data rows (2-76) of Table1 from
# "person-day" to "person-day/kW" on the "Tabelle1" worksheet.
# Column D holds the "unit" field.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -eq "person-day") {
        $cell.Value2 = "person-day/kW"
    }
}

$ws.Range("D9").Select()
